# Daily attendance processing - 2026-01-04 04:49:53
# Normalize the "Recorded By" (column G) text so the "System" author is
# listed first, ahead of the human reviewer's email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
